# atualizei dados bibi e add
# Replace the 4 rows for 2025-06-30 (rows 2-5) with the updated/refreshed
# dataset: the remaining rows are refreshed with new numbers and one new
# row (2025-07-15 / PROJETOR ASTRONAUTA HMASTON) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four obsolete "2025-06-30" rows entirely (rows 2-5).
$ws.Rows("2:5").Delete()

# Final data set (rows 2-8). Columns:
# A data | B quantidade_atipica | C cliente | D id_venda | E id_produto
# F produto | G estoque_atualizado | H media_vendas | I desvio_padrao
$rows = @(
    @("2025-07-01", 2, "BEMOL S/A",        "375040", 10114, "CARREGADOR USB-C A GOLD 20W CA31-4",                                      -90,  1.05, 0.22),
    @("2025-07-02", 2, "BEMOL S/A",        "375697", 13018, "ADAPTADOR TUDO EM UM UNIVERSAL INOVA PRIME TRA-30078",                     0,  1.08, 0.29),
    @("2025-07-07", 4, "BEMOL S/A",        "378212", 13546, "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON", -323,  1.1,  0.34),
    @("2025-07-09", 2, "BEMOL S/A",        "379513", 13000, "MOUSE PAD HARRY POTTER ESTAMPADO 26CMX21CM BLISTER C/1 UND LETRON",       -50,  1.03, 0.18),
    @("2025-07-09", 2, "MATHEUS SILVEIRA", "379106",  8915, "CAPA IPHONE 11",                                                           30,  1.06, 0.25),
    @("2025-07-11", 2, "BEMOL S/A",        "380683", 14241, "MOCHILA PELUCIA STITCH",                                                  -26,  1.04, 0.2),
    @("2025-07-15", 2, "BEMOL S/A",        "383049", 12016, "PROJETOR ASTRONAUTA HMASTON",                                             -40,  1.04, 0.21)
)

$r = 2
foreach ($row in $rows) {
    # Data and id_venda must stay text (not be reinterpreted as a date /
    # number), so use the quote-prefix trick, same as typing '... in Excel.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}
